# "error solve ifrs list" -- correct the IFRS financial figures in the
# 아세아제지 company_list sheet: rows 2-6 get corrected (much smaller,
# unit-corrected) financial figures, and the forecast rows 7-9 (which held
# erroneous duplicated/estimated data) have their financial columns cleared,
# leaving only the period labels in columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12): replace D2:AJ2 with corrected figures
$vals = @(6154,248,248,234,186,186,0,7157,2333,4824,4823,2,448,792,-643,-242,426,366,1398,4.03,3.02,3.91,2.6,48.37,976.83,2075,11.62,53995,0.45,450,1.87,21.63,8956502)
$arr = New-Object "object[,]" 1,$vals.Length
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("D2:AJ2").Value = $arr

# Row 3 (2015/12): replace D3:AJ3 with corrected figures
$vals = @(5977,-7,-7,-362,-387,-386,0,6546,2141,4404,4403,1,448,137,91,-464,190,-53,976,-0.11,-6.47,-8.380000000000001,-5.64,48.62,882.65,-4315,-3.99,49297,0.35,0,0,0,8956502)
$arr = New-Object "object[,]" 1,$vals.Length
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("D3:AJ3").Value = $arr

# Row 4 (2016/12): replace D4:AJ4 with corrected figures
$vals = @(6330,245,245,270,207,207,0,7210,2598,4611,4610,1,448,308,-433,504,598,-289,1493,3.87,3.28,4.6,3.02,56.34,929.24,2315,7.91,51616,0.35,500,2.73,21.54,8956502)
$arr = New-Object "object[,]" 1,$vals.Length
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("D4:AJ4").Value = $arr

# Row 5 (2017/12): replace D5:AJ5 with corrected figures
$vals = @(7286,54,54,130,50,50,0,7669,3061,4608,4607,1,448,-2,-471,446,628,-630,1971,0.74,0.6899999999999999,1.09,0.67,66.44,928.62,561,32.46,51433,0.35,450,2.47,80.25,8956502)
$arr = New-Object "object[,]" 1,$vals.Length
for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("D5:AJ5").Value = $arr

# Row 6 (2018/12): replace per-cell -- this row never had J6/O6 populated,
# and AG6/AH6 (현금DPS/현금배당수익률) are removed outright, not just zeroed.
$ws.Range("D6").Value = 7758
$ws.Range("E6").Value = 983
$ws.Range("F6").Value = 983
$ws.Range("G6").Value = 1031
$ws.Range("H6").Value = 787
$ws.Range("I6").Value = 787
$ws.Range("K6").Value = 8124
$ws.Range("L6").Value = 2789
$ws.Range("M6").Value = 5335
$ws.Range("N6").Value = 5335
$ws.Range("P6").Value = 448
$ws.Range("Q6").Value = 1462
$ws.Range("R6").Value = -618
$ws.Range("S6").Value = -452
$ws.Range("T6").Value = 602
$ws.Range("U6").Value = 860
$ws.Range("V6").Value = 1562
$ws.Range("W6").Value = 12.67
$ws.Range("X6").Value = 10.14
$ws.Range("Y6").Value = 15.83
$ws.Range("Z6").Value = 9.960000000000001
$ws.Range("AA6").Value = 52.26
$ws.Range("AB6").Value = 1091.01
$ws.Range("AC6").Value = 8785
$ws.Range("AD6").Value = 3.57
$ws.Range("AE6").Value = 59569
$ws.Range("AF6").Value = 0.53
$ws.Range("AI6").Value = 8.539999999999999
$ws.Range("AJ6").Value = 8956502
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E) estimates): clear all
# financial data columns, keeping only the A/B/C period labels.
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
